$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows for Cookie, Point and Ticket list models (mirroring the
# existing KingdomObjList row at row 23). Column A is filled in first for
# all three rows, then column B, then column E, to match the order new
# strings were appended to the shared-strings table.
$ws.Range("A24").Value = "CookieList"
$ws.Range("A25").Value = "PointList"
$ws.Range("A26").Value = "TicketList"

$ws.Range("B24").Value = "LIST:CookiePacket"
$ws.Range("B25").Value = "LIST:PointPacket"
$ws.Range("B26").Value = "LIST:TicketPacket"

$ws.Range("E24").Value = "Packet"
$ws.Range("E25").Value = "Packet"
$ws.Range("E26").Value = "Packet"

# Leave the selection on the last-edited cell, as captured in the workbook.
$ws.Range("E25").Select()
